$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.061.44"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "2.310.74"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'517.03"
$ws.Range("E5").Value = "  +4.19%  "
$ws.Range("D6").Value = "'133.06"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "2.331.67"
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("E10").Value = "  +8.56%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "'5.17"
$ws.Range("E12").Value = "  +7.61%  "
$ws.Range("D13").Value = "'0.342"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "'23.99"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "2.727.16"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "56.261.02"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("E17").Value = "  +4.37%  "
$ws.Range("D18").Value = "2.320.48"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").Value = "'10.53"
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("E20").Value = "  +3.52%  "
$ws.Range("D21").Value = "'321.16"
$ws.Range("E21").Value = "  +6.37%  "
$ws.Range("D22").Value = "'6.66"
$ws.Range("E22").Value = "  +5.00%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'0.992"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  +6.32%  "
$ws.Range("D27").Value = "'7.65"
$ws.Range("E27").Value = "  +4.29%  "
$ws.Range("D28").Value = "'171.23"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0732"
$ws.Range("E29").Value = "  +5.84%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.70"
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("D31").Value = "'1.19"
$ws.Range("E31").Value = "  +9.70%  "
$ws.Range("E32").Value = "  +5.08%  "
$ws.Range("D33").Value = "'18.32"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").Value = "'1.26"
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("D37").Value = "'0.925"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").Value = "'4.00"
$ws.Range("E38").Value = "  +7.89%  "
$ws.Range("E39").Value = "  +8.49%  "
$ws.Range("D40").Value = "'37.44"
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").Value = "'139.54"
$ws.Range("E42").Value = "  +11.41%  "
$ws.Range("D43").Value = "'3.57"
$ws.Range("E43").Value = "  +5.93%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'274.71"
$ws.Range("E44").Value = "  +13.53%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'5.09"
$ws.Range("E45").Value = "  +5.86%  "
$ws.Range("D46").Value = "'0.0509"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("D47").Value = "'0.0929"
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("D48").Value = "'0.556"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "'0.381"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").Value = "'0.0215"
$ws.Range("E50").Value = "  +5.02%  "
$ws.Range("E51").Value = "  +4.90%  "
